# Generate Report for Handback
# Marks the zh-cn and de-de localization rows as handed back: updates the
# Status text, fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns, links the target file back to the
# source markdown file, and widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shared by the Overview summary columns and both language sheets)
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6eac924673adab27a8c47bf1998e4fb00bfe6de/e2e/977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.md", "", "", "977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.md")
$ws2.Range("J2").Value = "977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.bd29fc494909e1a9accfe9744ca33797666fa5cd.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-29 12:59:17"

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f6eac924673adab27a8c47bf1998e4fb00bfe6de/e2e/977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.md", "", "", "977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.md")
$ws3.Range("J2").Value = "977574a9-6d7b-41ee-be0e-dfcb9cd2cb0a.bd29fc494909e1a9accfe9744ca33797666fa5cd.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-29 12:59:25"

# --- Widen columns that now contain the longer status text / file names ---
$ws1.Columns.Item(5).ColumnWidth = 29.17
$ws1.Columns.Item(6).ColumnWidth = 29.17
$ws2.Columns.Item(3).ColumnWidth = 29.17
$ws3.Columns.Item(3).ColumnWidth = 29.17
$ws2.Columns.Item(9).ColumnWidth = 39.17
$ws2.Columns.Item(10).ColumnWidth = 39.17
$ws3.Columns.Item(9).ColumnWidth = 39.17
$ws3.Columns.Item(10).ColumnWidth = 39.17
